# Update column F (dSF) values for rows where the data was repulled.
# Mapping of row -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 4
    8  = -1
    10 = -6
    12 = -10
    14 = -4
    15 = -5
    16 = 2
    17 = -3
    19 = -6
    20 = -2
    21 = -3
    24 = -4
    27 = 0
    31 = -6
    33 = -8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
